$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the reading form text in E6: "a", 123, etc… -> "a"
$ws.Range("E6").Value = '"a"'

# Move the active selection to E6
$ws.Range("E6").Select()
